$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# C10 holds the "From" value for rule R20; update it from 18 to 1.
$ws.Range("C10").Value = 1
